$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$c = $ws.Cells.Item(2, 4)
$c.NumberFormat = "@"
$c.Value = "43.912.35"
$c.Style = "Normal"
$c = $ws.Cells.Item(2, 5)
$c.NumberFormat = "@"
$c.Value = "  +0.11%  "
$c.Style = "Normal"
$c = $ws.Cells.Item(3, 4)
$c.NumberFormat = "@"
$c.Value = "2.271.16"
$c.Style = "Normal"
$c = $ws.Cells.Item(3, 5)
$c.NumberFormat = "@"
$c.Value = "  +2.49%  "
$c.Style = "Normal"
$c = $ws.Cells.Item(4, 5)
$c.NumberFormat = "@"
$c.Value = "  +0.14%  "
$c.Style = "Normal"
$c = $ws.Cells.Item(5, 4)
$c.NumberFormat = "@"
$c.Value = "268.33"
$c.Style = "Normal"
$c = $ws.Cells.Item(5, 5)
$c.NumberFormat = "@"
$c.Value = "  +2.61%  "
$c.Style = "Normal"
$c = $ws.Cells.Item(6, 4)
$c.NumberFormat = "@"
$c.Value = "93.45"
$c.Style = "Normal"
$c = $ws.Cells.Item(6, 5)
$c.NumberFormat = "@"
$c.Value = "  +7.96%  "
$c.Style = "Normal"
$c = $ws.Cells.Item(7, 4)
$c.NumberFormat = "@"
$c.Value = "0.625"
$c.Style = "Normal"
$c = $ws.Cells.Item(7, 5)
$c.NumberFormat = "@"
$c.Value = "  +0.76%  "
$c.Style = "Normal"
$c = $ws.Cells.Item(8, 5)
$c.NumberFormat = "@"
$c.Value = "  -0.01%  "
$c.Style = "Normal"
$c = $ws.Cells.Item(9, 4)
$c.NumberFormat = "@"
$c.Value = "0.619"
$c.Style = "Normal"
$c = $ws.Cells.Item(9, 5)
$c.NumberFormat = "@"
$c.Value = "  +2.05%  "
$c.Style = "Normal"
$c = $ws.Cells.Item(10, 4)
$c.NumberFormat = "@"
$c.Value = "46.49"
$c.Style = "Normal"
$c = $ws.Cells.Item(10, 5)
$c.NumberFormat = "@"
$c.Value = "  +2.86%  "
$c.Style = "Normal"
$c = $ws.Cells.Item(11, 4)
$c.NumberFormat = "@"
$c.Value = "0.0932"
$c.Style = "Normal"
$c = $ws.Cells.Item(11, 5)
$c.NumberFormat = "@"
$c.Value = "  +1.27%  "
$c.Style = "Normal"
$c = $ws.Cells.Item(12, 4)
$c.NumberFormat = "@"
$c.Value = "7.96"
$c.Style = "Normal"
$c = $ws.Cells.Item(12, 5)
$c.NumberFormat = "@"
$c.Value = "  +6.54%  "
$c.Style = "Normal"
$c = $ws.Cells.Item(13, 5)
$c.NumberFormat = "@"
$c.Value = "  +0.64%  "
$c.Style = "Normal"
$c = $ws.Cells.Item(14, 4)
$c.NumberFormat = "@"
$c.Value = "2.613.49"
$c.Style = "Normal"
$c = $ws.Cells.Item(14, 5)
$c.NumberFormat = "@"
$c.Value = "  +2.52%  "
$c.Style = "Normal"
$c = $ws.Cells.Item(15, 4)
$c.NumberFormat = "@"
$c.Value = "15.33"
$c.Style = "Normal"
$c = $ws.Cells.Item(15, 5)
$c.NumberFormat = "@"
$c.Value = "  +5.77%  "
$c.Style = "Normal"
$c = $ws.Cells.Item(16, 4)
$c.NumberFormat = "@"
$c.Value = "0.824"
$c.Style = "Normal"
$c = $ws.Cells.Item(16, 5)
$c.NumberFormat = "@"
$c.Value = "  +5.17%  "
$c.Style = "Normal"
$c = $ws.Cells.Item(17, 4)
$c.NumberFormat = "@"
$c.Value = "2.284.72"
$c.Style = "Normal"
$c = $ws.Cells.Item(17, 5)
$c.NumberFormat = "@"
$c.Value = "  +3.14%  "
$c.Style = "Normal"
$c = $ws.Cells.Item(18, 4)
$c.NumberFormat = "@"
$c.Value = "43.923.81"
$c.Style = "Normal"
$c = $ws.Cells.Item(18, 5)
$c.NumberFormat = "@"
$c.Value = "  +0.28%  "
$c.Style = "Normal"
$c = $ws.Cells.Item(19, 5)
$c.NumberFormat = "@"
$c.Value = "  +1.42%  "
$c.Style = "Normal"
$c = $ws.Cells.Item(20, 4)
$c.NumberFormat = "@"
$c.Value = "6.16"
$c.Style = "Normal"
$c = $ws.Cells.Item(20, 5)
$c.NumberFormat = "@"
$c.Value = "  +3.46%  "
$c.Style = "Normal"
$c = $ws.Cells.Item(21, 4)
$c.NumberFormat = "@"
$c.Value = "70.78"
$c.Style = "Normal"
$c = $ws.Cells.Item(21, 5)
$c.NumberFormat = "@"
$c.Value = "  +1.34%  "
$c.Style = "Normal"
$c = $ws.Cells.Item(22, 4)
$c.NumberFormat = "@"
$c.Value = "2.29"
$c.Style = "Normal"
$c = $ws.Cells.Item(22, 5)
$c.NumberFormat = "@"
$c.Value = "  -2.73%  "
$c.Style = "Normal"
$c = $ws.Cells.Item(23, 4)
$c.NumberFormat = "@"
$c.Value = "10.05"
$c.Style = "Normal"
$c = $ws.Cells.Item(23, 5)
$c.NumberFormat = "@"
$c.Value = "  +11.83%  "
$c.Style = "Normal"
$c = $ws.Cells.Item(24, 4)
$c.NumberFormat = "@"
$c.Value = "235.48"
$c.Style = "Normal"
$c = $ws.Cells.Item(24, 5)
$c.NumberFormat = "@"
$c.Value = "  +1.68%  "
$c.Style = "Normal"
$c = $ws.Cells.Item(25, 5)
$c.NumberFormat = "@"
$c.Value = "  -0.06%  "
$c.Style = "Normal"
$c = $ws.Cells.Item(26, 4)
$c.NumberFormat = "@"
$c.Value = "11.27"
$c.Style = "Normal"
$c = $ws.Cells.Item(26, 5)
$c.NumberFormat = "@"
$c.Value = "  +5.76%  "
$c.Style = "Normal"
$c = $ws.Cells.Item(27, 4)
$c.NumberFormat = "@"
$c.Value = "2.48"
$c.Style = "Normal"
$c = $ws.Cells.Item(27, 5)
$c.NumberFormat = "@"
$c.Value = "  +9.84%  "
$c.Style = "Normal"
$ws.Cells.Item(28, 2).Value = "WEMIXToken"
$ws.Cells.Item(28, 3).Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$c = $ws.Cells.Item(28, 4)
$c.NumberFormat = "@"
$c.Value = "3.35"
$c.Style = "Normal"
$c = $ws.Cells.Item(28, 5)
$c.NumberFormat = "@"
$c.Value = "  -5.41%  "
$c.Style = "Normal"
$ws.Cells.Item(29, 2).Value = "Toncoin"
$ws.Cells.Item(29, 3).Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$c = $ws.Cells.Item(29, 4)
$c.NumberFormat = "@"
$c.Value = "2.26"
$c.Style = "Normal"
$c = $ws.Cells.Item(29, 5)
$c.NumberFormat = "@"
$c.Value = "  +0.01%  "
$c.Style = "Normal"
$ws.Cells.Item(30, 2).Value = "InjectiveProtocol"
$ws.Cells.Item(30, 3).Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$c = $ws.Cells.Item(30, 4)
$c.NumberFormat = "@"
$c.Value = "39.05"
$c.Style = "Normal"
$c = $ws.Cells.Item(30, 5)
$c.NumberFormat = "@"
$c.Value = "  -3.17%  "
$c.Style = "Normal"
$ws.Cells.Item(31, 2).Value = "Monero"
$ws.Cells.Item(31, 3).Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$c = $ws.Cells.Item(31, 4)
$c.NumberFormat = "@"
$c.Value = "173.78"
$c.Style = "Normal"
$c = $ws.Cells.Item(31, 5)
$c.NumberFormat = "@"
$c.Value = "  -0.34%  "
$c.Style = "Normal"
$ws.Cells.Item(32, 2).Value = "EthereumClassic"
$ws.Cells.Item(32, 3).Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$c = $ws.Cells.Item(32, 4)
$c.NumberFormat = "@"
$c.Value = "22.00"
$c.Style = "Normal"
$c = $ws.Cells.Item(32, 5)
$c.NumberFormat = "@"
$c.Value = "  +7.40%  "
$c.Style = "Normal"
$c = $ws.Cells.Item(33, 4)
$c.NumberFormat = "@"
$c.Value = "0.0904"
$c.Style = "Normal"
$c = $ws.Cells.Item(33, 5)
$c.NumberFormat = "@"
$c.Value = "  +4.03%  "
$c.Style = "Normal"
$c = $ws.Cells.Item(34, 4)
$c.NumberFormat = "@"
$c.Value = "5.57"
$c.Style = "Normal"
$c = $ws.Cells.Item(34, 5)
$c.NumberFormat = "@"
$c.Value = "  +3.15%  "
$c.Style = "Normal"
$c = $ws.Cells.Item(35, 4)
$c.NumberFormat = "@"
$c.Value = "0.124"
$c.Style = "Normal"
$c = $ws.Cells.Item(35, 5)
$c.NumberFormat = "@"
$c.Value = "  +1.01%  "
$c.Style = "Normal"
$c = $ws.Cells.Item(36, 4)
$c.NumberFormat = "@"
$c.Value = "0.111"
$c.Style = "Normal"
$c = $ws.Cells.Item(36, 5)
$c.NumberFormat = "@"
$c.Value = "  -1.18%  "
$c.Style = "Normal"
$c = $ws.Cells.Item(37, 4)
$c.NumberFormat = "@"
$c.Value = "4.41"
$c.Style = "Normal"
$c = $ws.Cells.Item(37, 5)
$c.NumberFormat = "@"
$c.Value = "  -1.89%  "
$c.Style = "Normal"
$c = $ws.Cells.Item(38, 4)
$c.NumberFormat = "@"
$c.Value = "0.0349"
$c.Style = "Normal"
$c = $ws.Cells.Item(38, 5)
$c.NumberFormat = "@"
$c.Value = "  -2.05%  "
$c.Style = "Normal"
$c = $ws.Cells.Item(39, 4)
$c.NumberFormat = "@"
$c.Value = "3.43"
$c.Style = "Normal"
$c = $ws.Cells.Item(39, 5)
$c.NumberFormat = "@"
$c.Value = "  +17.22%  "
$c.Style = "Normal"
$c = $ws.Cells.Item(40, 4)
$c.NumberFormat = "@"
$c.Value = "0.246"
$c.Style = "Normal"
$c = $ws.Cells.Item(40, 5)
$c.NumberFormat = "@"
$c.Value = "  +22.46%  "
$c.Style = "Normal"
$c = $ws.Cells.Item(41, 4)
$c.NumberFormat = "@"
$c.Value = "2.22"
$c.Style = "Normal"
$c = $ws.Cells.Item(41, 5)
$c.NumberFormat = "@"
$c.Value = "  +6.06%  "
$c.Style = "Normal"
$c = $ws.Cells.Item(42, 4)
$c.NumberFormat = "@"
$c.Value = "12.30"
$c.Style = "Normal"
$c = $ws.Cells.Item(42, 5)
$c.NumberFormat = "@"
$c.Value = "  -3.02%  "
$c.Style = "Normal"
$c = $ws.Cells.Item(43, 4)
$c.NumberFormat = "@"
$c.Value = "5.46"
$c.Style = "Normal"
$c = $ws.Cells.Item(43, 5)
$c.NumberFormat = "@"
$c.Value = "  -1.15%  "
$c.Style = "Normal"
$c = $ws.Cells.Item(44, 4)
$c.NumberFormat = "@"
$c.Value = "61.08"
$c.Style = "Normal"
$c = $ws.Cells.Item(44, 5)
$c.NumberFormat = "@"
$c.Value = "  -3.51%  "
$c.Style = "Normal"
$ws.Cells.Item(45, 2).Value = "ARBITRUM"
$ws.Cells.Item(45, 3).Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$c = $ws.Cells.Item(45, 4)
$c.NumberFormat = "@"
$c.Value = "1.22"
$c.Style = "Normal"
$c = $ws.Cells.Item(45, 5)
$c.NumberFormat = "@"
$c.Value = "  +8.72%  "
$c.Style = "Normal"
$ws.Cells.Item(46, 2).Value = "Cronos"
$ws.Cells.Item(46, 3).Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$c = $ws.Cells.Item(46, 4)
$c.NumberFormat = "@"
$c.Value = "0.102"
$c.Style = "Normal"
$c = $ws.Cells.Item(46, 5)
$c.NumberFormat = "@"
$c.Value = "  +4.07%  "
$c.Style = "Normal"
$c = $ws.Cells.Item(47, 4)
$c.NumberFormat = "@"
$c.Value = "8.52"
$c.Style = "Normal"
$c = $ws.Cells.Item(47, 5)
$c.NumberFormat = "@"
$c.Value = "  +2.10%  "
$c.Style = "Normal"
$c = $ws.Cells.Item(48, 4)
$c.NumberFormat = "@"
$c.Value = "99.19"
$c.Style = "Normal"
$c = $ws.Cells.Item(48, 5)
$c.NumberFormat = "@"
$c.Value = "  -1.53%  "
$c.Style = "Normal"
$c = $ws.Cells.Item(49, 4)
$c.NumberFormat = "@"
$c.Value = "1.18"
$c.Style = "Normal"
$c = $ws.Cells.Item(49, 5)
$c.NumberFormat = "@"
$c.Value = "  -0.17%  "
$c.Style = "Normal"
$c = $ws.Cells.Item(50, 4)
$c.NumberFormat = "@"
$c.Value = "0.427"
$c.Style = "Normal"
$c = $ws.Cells.Item(50, 5)
$c.NumberFormat = "@"
$c.Value = "  -2.87%  "
$c.Style = "Normal"
$c = $ws.Cells.Item(51, 4)
$c.NumberFormat = "@"
$c.Value = "2.493.65"
$c.Style = "Normal"
$c = $ws.Cells.Item(51, 5)
$c.NumberFormat = "@"
$c.Value = "  +2.42%  "
$c.Style = "Normal"
